# Update column Q ('uds. Objetivo semana pasada') values per the commit diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("Q3").Value = 1
$ws.Range("Q6").Value = 6
$ws.Range("Q7").Value = 1
$ws.Range("Q15").Value = 1
$ws.Range("Q21").Value = 1
$ws.Range("Q28").Value = 1
$ws.Range("Q29").Value = 1
$ws.Range("Q32").Value = 1
$ws.Range("Q38").Value = 1
$ws.Range("Q41").Value = 2
$ws.Range("Q43").Value = 1
$ws.Range("Q44").Value = 1
$ws.Range("Q48").Value = 6
$ws.Range("Q49").Value = 3
$ws.Range("Q50").Value = 2
$ws.Range("Q55").Value = 2
$ws.Range("Q62").Value = 1
$ws.Range("Q63").Value = 1
$ws.Range("Q69").Value = 1
$ws.Range("Q73").Value = 1
$ws.Range("Q74").Value = 5
$ws.Range("Q83").Value = 5
$ws.Range("Q84").Value = 1
$ws.Range("Q86").Value = 1
$ws.Range("Q100").Value = 1
$ws.Range("Q106").Value = 3
$ws.Range("Q110").Value = 1
$ws.Range("Q112").Value = 1
$ws.Range("Q113").Value = 1
$ws.Range("Q117").Value = 1
$ws.Range("Q118").Value = 1
$ws.Range("Q121").Value = 1
$ws.Range("Q123").Value = 1
$ws.Range("Q125").Value = 1
$ws.Range("Q126").Value = 3
$ws.Range("Q127").Value = 3
$ws.Range("Q128").Value = 2
$ws.Range("Q129").Value = 3
$ws.Range("Q130").Value = 1
$ws.Range("Q131").Value = 5
$ws.Range("Q134").Value = 1
$ws.Range("Q136").Value = 1
$ws.Range("Q137").Value = 19
$ws.Range("Q138").Value = 3
$ws.Range("Q139").Value = 1
$ws.Range("Q141").Value = 1
$ws.Range("Q142").Value = 5
$ws.Range("Q143").Value = 3
$ws.Range("Q144").Value = 1
$ws.Range("Q148").Value = 2
$ws.Range("Q151").Value = 1
$ws.Range("Q155").Value = 1
$ws.Range("Q157").Value = 1
$ws.Range("Q164").Value = 1
$ws.Range("Q165").Value = 2
$ws.Range("Q169").Value = 3
$ws.Range("Q171").Value = 2
$ws.Range("Q176").Value = 1
$ws.Range("Q178").Value = 1
$ws.Range("Q180").Value = 1
$ws.Range("Q181").Value = 1
$ws.Range("Q193").Value = 1
$ws.Range("Q194").Value = 1
$ws.Range("Q198").Value = 1
$ws.Range("Q199").Value = 2
$ws.Range("Q200").Value = 3
$ws.Range("Q201").Value = 1
$ws.Range("Q202").Value = 2
$ws.Range("Q204").Value = 2
$ws.Range("Q205").Value = 1
$ws.Range("Q207").Value = 2
$ws.Range("Q209").Value = 1
$ws.Range("Q210").Value = 1
$ws.Range("Q211").Value = 7
$ws.Range("Q212").Value = 1
$ws.Range("Q215").Value = 1
$ws.Range("Q216").Value = 1
$ws.Range("Q219").Value = 1
$ws.Range("Q220").Value = 1
